# Apply the "more error handling" update to the Orders sheet:
# - A search for a product that doesn't exist ("Ipoh Coff" instead of
#   "Ipoh Coffee") now shows up as a Failed order, while every other
#   order (previously a mix of Succeeded/Failed) now shows Succeeded.
# - The workbook window is left minimized and the last selection is G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Row 3: change status from Failed to Succeeded
$ws.Range("C3").Value = "Succeeded"

# Rows 4-8: add a Succeeded status (previously blank)
$ws.Range("C4").Value = "Succeeded"
$ws.Range("C5").Value = "Succeeded"
$ws.Range("C6").Value = "Succeeded"
$ws.Range("C7").Value = "Succeeded"
$ws.Range("C8").Value = "Succeeded"

# Row 9: product name becomes a non-existent product ("Ipoh Coff"),
# which fails the order lookup, so the status becomes Failed.
$ws.Range("A9").Value = "Ipoh Coff"
$ws.Range("C9").Value = "Failed"

# Update the last active selection on the Orders sheet to G6
$ws.Activate()
[void]$ws.Range("G6").Select()

# Minimize the workbook window
$win = $wb.Windows.Item(1)
$win.WindowState = -4140
